# Updated cryptos list values (Price / Volume(1h)) per target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.352.39"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "2.314.93"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.43%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.974"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.54%  "

$ws.Range("E15").Value = "  -4.65%  "

$ws.Range("D16").Value = "2.666.92"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").Value = "2.348.94"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "42.310.84"
$ws.Range("E18").Value = "  -0.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.39%  "

$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "261.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.97%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0896"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("E33").Value = "  -6.44%  "

$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.119"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.35%  "

$ws.Range("E36").Value = "  -3.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.22%  "

$ws.Range("E40").Value = "  -10.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.31%  "

$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.232"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "111.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "
